$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'256.22"
$ws.Cells.Item(2, 5).Value = "'4.51%"
$ws.Cells.Item(3, 4).Value = "'28.05"
$ws.Cells.Item(3, 5).Value = "'-3.48%"
$ws.Cells.Item(4, 4).Value = "'5.218"
$ws.Cells.Item(4, 5).Value = "'-0.96%"
$ws.Cells.Item(5, 4).Value = "'0.05887"
$ws.Cells.Item(5, 5).Value = "'3.05%"
$ws.Cells.Item(6, 4).Value = "'6.694"
$ws.Cells.Item(6, 5).Value = "'1.35%"
$ws.Cells.Item(7, 4).Value = "'0.8695"
$ws.Cells.Item(7, 5).Value = "'1.90%"
$ws.Cells.Item(8, 4).Value = "'0.9743"
$ws.Cells.Item(8, 5).Value = "'13.56%"
$ws.Cells.Item(9, 4).Value = "'0.1409"
$ws.Cells.Item(9, 5).Value = "'2.67%"
$ws.Cells.Item(10, 4).Value = "'0.07147"
$ws.Cells.Item(10, 5).Value = "'1.07%"
$ws.Cells.Item(11, 4).Value = "'0.03174"
$ws.Cells.Item(11, 5).Value = "'0.20%"
$ws.Cells.Item(12, 4).Value = "'0.09222"
$ws.Cells.Item(12, 5).Value = "'-0.70%"
$ws.Cells.Item(13, 4).Value = "'0.001537"
$ws.Cells.Item(13, 5).Value = "'0.80%"
$ws.Cells.Item(14, 4).Value = "'0.0006084"
$ws.Cells.Item(14, 5).Value = "'1.60%"
$ws.Cells.Item(15, 4).Value = "'0.006026"
$ws.Cells.Item(15, 5).Value = "'1.10%"
$ws.Cells.Item(16, 4).Value = "'3.498"
$ws.Cells.Item(16, 5).Value = "'-0.30%"
$ws.Cells.Item(17, 4).Value = "'3.216"
$ws.Cells.Item(17, 5).Value = "'1.31%"
$ws.Cells.Item(18, 5).Value = "'2.03%"
$ws.Cells.Item(19, 4).Value = "'0.3172"
$ws.Cells.Item(19, 5).Value = "'0.03%"
$ws.Cells.Item(20, 4).Value = "'0.03490"
$ws.Cells.Item(20, 5).Value = "'5.03%"
$ws.Cells.Item(21, 4).Value = "'0.1289"
$ws.Cells.Item(21, 5).Value = "'0.95%"
$ws.Cells.Item(22, 4).Value = "'3.536"
$ws.Cells.Item(22, 5).Value = "'1.42%"
$ws.Cells.Item(23, 4).Value = "'0.04170"
$ws.Cells.Item(23, 5).Value = "'1.22%"
$ws.Cells.Item(24, 4).Value = "'0.1365"
$ws.Cells.Item(24, 5).Value = "'2.71%"
$ws.Cells.Item(25, 4).Value = "'0.001226"
$ws.Cells.Item(25, 5).Value = "'0.57%"
$ws.Cells.Item(26, 4).Value = "'0.004558"
$ws.Cells.Item(26, 5).Value = "'9.92%"
$ws.Cells.Item(27, 4).Value = "'0.0001199"
$ws.Cells.Item(27, 5).Value = "'-0.05%"
$ws.Cells.Item(28, 4).Value = "'0.0001471"
$ws.Cells.Item(28, 5).Value = "'1.55%"
$ws.Cells.Item(40, 4).Value = "'0.03829"
$ws.Cells.Item(40, 5).Value = "'0.65%"
$ws.Cells.Item(41, 4).Value = "'0.005475"
$ws.Cells.Item(41, 5).Value = "'-4.24%"
$ws.Cells.Item(42, 5).Value = "'3.63%"
$ws.Cells.Item(43, 4).Value = "'0.002342"
$ws.Cells.Item(43, 5).Value = "'-3.02%"
$ws.Cells.Item(44, 4).Value = "'0.009558"
$ws.Cells.Item(44, 5).Value = "'1.80%"
$ws.Cells.Item(45, 4).Value = "'0.00005386"
$ws.Cells.Item(45, 5).Value = "'2.11%"
$ws.Cells.Item(46, 4).Value = "'0.00000000749"
$ws.Cells.Item(46, 5).Value = "'-0.06%"
$ws.Cells.Item(47, 4).Value = "'0.09491"
$ws.Cells.Item(47, 5).Value = "'5.61%"
$ws.Cells.Item(48, 4).Value = "'0.002129"
$ws.Cells.Item(49, 4).Value = "'0.00002098"
$ws.Cells.Item(49, 5).Value = "'-0.06%"
$ws.Cells.Item(50, 4).Value = "'0.0001998"
$ws.Cells.Item(50, 5).Value = "'-0.06%"
